{"js": "// Update the date line and the 25 \"dividend\u00f7divisor=quotient, remainder\"\n// answer cells in the table. Every \"old\" string below is unique within the\n// document, so a plain text search-and-replace (no wildcards) is safe and\n// unambiguous.\nconst replacements = [\n  [\"2025-11-27 Thursday\", \"2025-11-28 Friday\"],\n  [\"342\u00f77=48, 6\", \"975\u00f72=487, 1\"],\n  [\"849\u00f77=121, 2\", \"779\u00f77=111, 2\"],\n  [\"682\u00f73=227, 1\", \"662\u00f77=94, 4\"],\n  [\"367\u00f72=183, 1\", \"909\u00f72=454, 1\"],\n  [\"576\u00f77=82, 2\", \"117\u00f77=16, 5\"],\n  [\"570\u00f74=142, 2\", \"176\u00f73=58, 2\"],\n  [\"917\u00f75=183, 2\", \"454\u00f72=227, 0\"],\n  [\"575\u00f75=115, 0\", \"966\u00f79=107, 3\"],\n  [\"188\u00f79=20, 8\", \"859\u00f76=143, 1\"],\n  [\"466\u00f72=233, 0\", \"166\u00f79=18, 4\"],\n  [\"847\u00f77=121, 0\", \"534\u00f76=89, 0\"],\n  [\"507\u00f77=72, 3\", \"606\u00f76=101, 0\"],\n  [\"920\u00f72=460, 0\", \"534\u00f77=76, 2\"],\n  [\"421\u00f73=140, 1\", \"791\u00f79=87, 8\"],\n  [\"804\u00f72=402, 0\", \"776\u00f79=86, 2\"],\n  [\"888\u00f74=222, 0\", \"885\u00f75=177, 0\"],\n  [\"711\u00f78=88, 7\", \"346\u00f73=115, 1\"],\n  [\"396\u00f78=49, 4\", \"651\u00f72=325, 1\"],\n  [\"115\u00f78=14, 3\", \"343\u00f73=114, 1\"],\n  [\"538\u00f72=269, 0\", \"259\u00f73=86, 1\"],\n  [\"784\u00f79=87, 1\", \"112\u00f76=18, 4\"],\n  [\"608\u00f76=101, 2\", \"411\u00f73=137, 0\"],\n  [\"399\u00f75=79, 4\", \"719\u00f78=89, 7\"],\n  [\"442\u00f75=88, 2\", \"266\u00f78=33, 2\"],\n  [\"431\u00f72=215, 1\", \"533\u00f75=106, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 \"dividend\u00f7divisor=quotient, remainder\"\n# answer cells in the table. Every \"old\" string is unique within the\n# document, so Find/Replace (no wildcards, whole document) is safe and\n# unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-11-27 Thursday\", \"2025-11-28 Friday\"),\n    @(\"342\u00f77=48, 6\", \"975\u00f72=487, 1\"),\n    @(\"849\u00f77=121, 2\", \"779\u00f77=111, 2\"),\n    @(\"682\u00f73=227, 1\", \"662\u00f77=94, 4\"),\n    @(\"367\u00f72=183, 1\", \"909\u00f72=454, 1\"),\n    @(\"576\u00f77=82, 2\", \"117\u00f77=16, 5\"),\n    @(\"570\u00f74=142, 2\", \"176\u00f73=58, 2\"),\n    @(\"917\u00f75=183, 2\", \"454\u00f72=227, 0\"),\n    @(\"575\u00f75=115, 0\", \"966\u00f79=107, 3\"),\n    @(\"188\u00f79=20, 8\", \"859\u00f76=143, 1\"),\n    @(\"466\u00f72=233, 0\", \"166\u00f79=18, 4\"),\n    @(\"847\u00f77=121, 0\", \"534\u00f76=89, 0\"),\n    @(\"507\u00f77=72, 3\", \"606\u00f76=101, 0\"),\n    @(\"920\u00f72=460, 0\", \"534\u00f77=76, 2\"),\n    @(\"421\u00f73=140, 1\", \"791\u00f79=87, 8\"),\n    @(\"804\u00f72=402, 0\", \"776\u00f79=86, 2\"),\n    @(\"888\u00f74=222, 0\", \"885\u00f75=177, 0\"),\n    @(\"711\u00f78=88, 7\", \"346\u00f73=115, 1\"),\n    @(\"396\u00f78=49, 4\", \"651\u00f72=325, 1\"),\n    @(\"115\u00f78=14, 3\", \"343\u00f73=114, 1\"),\n    @(\"538\u00f72=269, 0\", \"259\u00f73=86, 1\"),\n    @(\"784\u00f79=87, 1\", \"112\u00f76=18, 4\"),\n    @(\"608\u00f76=101, 2\", \"411\u00f73=137, 0\"),\n    @(\"399\u00f75=79, 4\", \"719\u00f78=89, 7\"),\n    @(\"442\u00f75=88, 2\", \"266\u00f78=33, 2\"),\n    @(\"431\u00f72=215, 1\", \"533\u00f75=106, 3\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
